$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Anzahl der Stunden mit Flaute" column for the first turbine (column K).
# This shifts the second turbine block (L:S) one column to the left (K:R).
$ws.Columns("K").Delete()

# Remove the "Anzahl der Stunden mit Flaute" column for the second turbine, which is
# now column R after the previous shift.
$ws.Columns("R").Delete()

# Update the recalculated energy yield figures for the first turbine (column H = E).
$ws.Range("H4").Value = 0.00288333673660529
$ws.Range("H5").Value = 0.05975693179000923
$ws.Range("H6").Value = 0.204366939880618
$ws.Range("H7").Value = 0.4211092519908928
$ws.Range("H8").Value = 0.6620772456870806
$ws.Range("H9").Value = 0.8472720869607887
$ws.Range("H10").Value = 0.9041900162899342
$ws.Range("H11").Value = 0.8045686718093682
$ws.Range("H12").Value = 0.6075227158346065
$ws.Range("H13").Value = 0.388342090392869
$ws.Range("H14").Value = 0.2127745694336207
$ws.Range("H15").Value = 0.1050963753924294
$ws.Range("H16").Value = 0.04846082529813664
$ws.Range("H17").Value = 0.0208778073951871
$ws.Range("H18").Value = 0.00840924374147051
$ws.Range("H19").Value = 0.003168412241897542
$ws.Range("H20").Value = 0.001117206728181465
$ws.Range("H21").Value = 0.0003688038001986859
$ws.Range("H22").Value = 0.0001140157908857068
$ws.Range("H23").Value = [double]"3.301870210998251e-05"
$ws.Range("H24").Value = [double]"8.95945416576852e-06"
$ws.Range("H25").Value = [double]"2.27833377323346e-06"
$ws.Range("H26").Value = [double]"5.430550634626768e-07"
$ws.Range("H27").Value = [double]"1.213467700566818e-07"

# Update the recalculated hourly energy figures for the first turbine (column I = E_h).
$ws.Range("I13").Value = 0.07323725000079946
$ws.Range("I14").Value = 0.04012705478218409
$ws.Range("I16").Value = 0.009139203978673004
$ws.Range("I17").Value = 0.003937335760136122
$ws.Range("I19").Value = 0.0005975293567346588
$ws.Range("I20").Value = 0.0002106934851476599
$ws.Range("I21").Value = [double]"6.95525331520748e-05"
$ws.Range("I23").Value = [double]"6.226981316097685e-06"

